$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to text format before writing, so that
# values like "0.9980" or "1.000" are stored as text (matching the
# original inlineStr / shared-string cells) rather than being
# auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.627.46"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").Value = "1.863.54"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "234.53"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "0.4691"
$ws.Range("E7").Value = "  -1.59%  "

$ws.Range("D8").Value = "0.2758"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").Value = "0.06359"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  +8.40%  "

$ws.Range("D11").Value = "1.852.63"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").Value = "0.07461"
$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").Value = "4.960"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").Value = "85.00"
$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("D15").Value = "0.6309"
$ws.Range("E15").Value = "  -0.24%  "

$ws.Range("D16").Value = "30.561.96"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").Value = "242.28"
$ws.Range("E17").Value = "  +4.46%  "

$ws.Range("D18").Value = "0.9979"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "0.000007364"
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "4.986"
$ws.Range("E22").Value = "  -2.40%  "

$ws.Range("D23").Value = "5.963"
$ws.Range("E23").Value = "  -1.00%  "

$ws.Range("D24").Value = "9.266"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("D25").Value = "166.87"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").Value = "18.19"
$ws.Range("E26").Value = "  +1.63%  "

$ws.Range("D27").Value = "1.887"
$ws.Range("E27").Value = "  +1.52%  "

$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("D29").Value = "1.374"
$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("D30").Value = "4.112"
$ws.Range("E30").Value = "  -2.99%  "

$ws.Range("D31").Value = "3.861"
$ws.Range("E31").Value = "  -1.39%  "

$ws.Range("D32").Value = "0.04931"
$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").Value = "1.151"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "0.7106"
$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("D35").Value = "2.698"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "0.01913"
$ws.Range("E36").Value = "  -2.47%  "

$ws.Range("D37").Value = "2.693"
$ws.Range("E37").Value = "  +2.40%  "

$ws.Range("D38").Value = "0.8829"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("D39").Value = "1.981"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").Value = "105.80"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").Value = "0.9980"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").Value = "0.4098"
$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").Value = "7.201"
$ws.Range("E44").Value = "  +1.86%  "

$ws.Range("D45").Value = "0.1238"
$ws.Range("E45").Value = "  +2.45%  "

$ws.Range("D46").Value = "61.98"
$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("D47").Value = "33.71"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("E48").Value = "  -2.48%  "

$ws.Range("D49").Value = "0.05556"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").Value = "1.376"
$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("D51").Value = "0.3708"
$ws.Range("E51").Value = "  -0.18%  "

# Restore the default (Normal) style so no explicit cell style is left
# behind, matching the unstyled D/E data cells in the original file.
$dataRange.Style = "Normal"
